$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Coin name / Link columns (B, C) for the rows whose coin ---
# --- order/identity changed between runs (plain text, safe as-is) ---
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("B47").Value = 'Monero'
$ws.Range("C47").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("B50").Value = 'ImmutableX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("B51").Value = 'FLOKI'
$ws.Range("C51").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'

# --- Update Price / Volume(1h) columns (D, E). These often look like ---
# --- numbers ("1.00", "0.999", ...) so force text format first to  ---
# --- preserve them as literal strings (matching the source data),  ---
# --- then restore the default "Normal" style so no stray number    ---
# --- format is left behind on the cell.                            ---
$priceVolRange = $ws.Range("D2:E51")
$priceVolRange.NumberFormat = "@"

$ws.Range("D2").Value = '90.622.14'
$ws.Range("E2").Value = '  +2.47%  '
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = '213.37'
$ws.Range("E5").Value = '  +0.13%  '
$ws.Range("D6").Value = '614.65'
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("D7").Value = '0.372'
$ws.Range("E7").Value = '  -4.45%  '
$ws.Range("D8").Value = '0.898'
$ws.Range("E8").Value = '  +15.72%  '
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.12%  '
$ws.Range("D10").Value = '3.040.23'
$ws.Range("E10").Value = '  -2.31%  '
$ws.Range("D11").Value = '0.678'
$ws.Range("E11").Value = '  +21.22%  '
$ws.Range("E12").Value = '  +5.67%  '
$ws.Range("D13").Value = '0.0000244'
$ws.Range("E13").Value = '  -2.03%  '
$ws.Range("D14").Value = '5.36'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '90.291.84'
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = '32.59'
$ws.Range("E16").Value = '  +1.90%  '
$ws.Range("D17").Value = '3.607.93'
$ws.Range("E17").Value = '  -2.28%  '
$ws.Range("D18").Value = '2.990.92'
$ws.Range("E18").Value = '  -4.32%  '
$ws.Range("D19").Value = '3.38'
$ws.Range("E19").Value = '  -0.99%  '
$ws.Range("D20").Value = '0.0000220'
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").Value = '13.55'
$ws.Range("E21").Value = '  +3.23%  '
$ws.Range("D22").Value = '428.37'
$ws.Range("E22").Value = '  +1.97%  '
$ws.Range("D23").Value = '8.34'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = '5.07'
$ws.Range("E24").Value = '  +3.81%  '
$ws.Range("D25").Value = '5.41'
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("D26").Value = '83.37'
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").Value = '11.71'
$ws.Range("E27").Value = '  +2.84%  '
$ws.Range("E28").Value = '  -2.76%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("D30").Value = '0.166'
$ws.Range("E30").Value = '  +6.73%  '
$ws.Range("E31").Value = '  +0.24%  '
$ws.Range("D32").Value = '8.69'
$ws.Range("E32").Value = '  +7.03%  '
$ws.Range("D33").Value = '3.77'
$ws.Range("E33").Value = '  -5.33%  '
$ws.Range("D34").Value = '506.60'
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = '6.77'
$ws.Range("E35").Value = '  -1.51%  '
$ws.Range("D36").Value = '22.96'
$ws.Range("E36").Value = '  +3.50%  '
$ws.Range("D37").Value = '1.82'
$ws.Range("E37").Value = '  -0.35%  '
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").Value = '0.134'
$ws.Range("E39").Value = '  -8.83%  '
$ws.Range("D40").Value = '22.31'
$ws.Range("E40").Value = '  +0.43%  '
$ws.Range("D41").Value = '0.999'
$ws.Range("E41").Value = '  -0.31%  '
$ws.Range("E42").Value = '  +0.00%  '
$ws.Range("D43").Value = '0.139'
$ws.Range("E43").Value = '  +5.88%  '
$ws.Range("D44").Value = '0.362'
$ws.Range("E44").Value = '  -0.15%  '
$ws.Range("D45").Value = '1.85'
$ws.Range("E45").Value = '  +0.51%  '
$ws.Range("D46").Value = '0.0708'
$ws.Range("E46").Value = '  +7.27%  '
$ws.Range("D47").Value = '143.70'
$ws.Range("E47").Value = '  -1.48%  '
$ws.Range("D48").Value = '43.57'
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").Value = '4.24'
$ws.Range("E49").Value = '  +8.30%  '
$ws.Range("D50").Value = '1.23'
$ws.Range("E50").Value = '  +3.87%  '
$ws.Range("D51").Value = '0.000253'
$ws.Range("E51").Value = '  +6.83%  '

$priceVolRange.Style = "Normal"

